# "Chuc nang chinh ta gan hoan thien" -- fill in the previously-blank
# "ĐG lần 3 (13/05)" (column H) scores for the 5 students in rows 11-15 of
# the "Ghi cong" sheet, and move the current selection to reflect where the
# user left off working (I15).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ghi cong")
$ws.Activate()

$ws.Range("H11").Value = 0
$ws.Range("H12").Value = 0
$ws.Range("H13").Value = 0
$ws.Range("H14").Value = 1
$ws.Range("H15").Value = 0

$ws.Range("I15").Select()
